$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) The paragraph about Duna-ag bridges used to be split into two runs
#    around a stray "_GoBack" bookmark, breaking the word "segito" into
#    "se" | "gito". Replacing across both runs re-joins them into one
#    contiguous run with the word spelled correctly, and removes the
#    bookmark that used to sit between them.
# ---------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.Execute("közlekedést segítő legfontosabb", $true, $false, $false, $false, $false, $true, 1, $false, "közlekedést segítő legfontosabb", 2)

# ---------------------------------------------------------------------
# 2) Nudge the floating picture further down the page (vertical offset
#    goes from 1654695 EMU to 2164715 EMU).
# ---------------------------------------------------------------------
foreach ($shp in $d.Shapes) {
    if ($shp.Name -eq "Kép 1") {
        $shp.Top = 2164715 / 914400 * 72
    }
}

# ---------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark later in the document, right after
#    the "E" of "Elsonek" in the Lanchid paragraph. Word only ever keeps
#    one "_GoBack" bookmark, so adding it here both relocates it and
#    splits the enclosing run into "...kepviselok. E" | "lsonek...".
# ---------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.Execute("Elsőnek", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $find2.Parent.Start + 1
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))
